$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-12 Saturday" "2025-04-13 Sunday"

Replace-Text "92×80=" "73×86="
Replace-Text "19×99=" "49×97="
Replace-Text "31×77=" "84×92="
Replace-Text "45×24=" "76×36="
Replace-Text "78×31=" "19×35="
Replace-Text "13×50=" "42×47="
Replace-Text "83×97=" "74×53="
Replace-Text "86×79=" "71×63="
Replace-Text "88×62=" "99×95="
Replace-Text "64×55=" "53×30="
Replace-Text "79×24=" "51×93="
Replace-Text "24×94=" "40×12="
Replace-Text "89×89=" "94×45="
Replace-Text "20×47=" "67×85="
Replace-Text "18×71=" "39×66="
Replace-Text "13×78=" "46×42="
Replace-Text "98×14=" "53×32="
Replace-Text "64×37=" "11×91="
Replace-Text "36×59=" "43×25="
Replace-Text "34×81=" "36×63="
Replace-Text "67×43=" "56×76="
Replace-Text "45×70=" "40×46="
Replace-Text "91×37=" "72×91="
Replace-Text "46×88=" "72×49="
Replace-Text "31×11=" "11×97="
